$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Url img" values (column F) for rows 2..21, pointing at the new
# GitHub-hosted images instead of the old caelum-online-public S3 bucket.
$urls = @{
    2  = "https://github.com/rosanamassahud/cursoPowerBI_base1/blob/main/img_petshop/Ra%C3%A7%C3%A3o%2Bfilhote%2BGatito.png"
    3  = "https://github.com/rosanamassahud/cursoPowerBI_base1/blob/main/img_petshop/Ra%C3%A7%C3%A3o%2Badulto%2BGatito.png"
    4  = "https://github.com/rosanamassahud/cursoPowerBI_base1/blob/main/img_petshop/Ra%C3%A7%C3%A3o%2Bs%C3%AAnior%2BGatito.png"
    5  = "https://github.com/rosanamassahud/cursoPowerBI_base1/blob/main/img_petshop/Petisco%2BGatito.png"
    6  = "https://github.com/rosanamassahud/cursoPowerBI_base1/blob/main/img_petshop/Rato%2BGatito.png"
    7  = "https://github.com/rosanamassahud/cursoPowerBI_base1/blob/main/img_petshop/Varinha%2BGatito.png"
    8  = "https://github.com/rosanamassahud/cursoPowerBI_base1/blob/main/img_petshop/Roupa%2Bde%2Binverno%2BGatito.png"
    9  = "https://github.com/rosanamassahud/cursoPowerBI_base1/blob/main/img_petshop/Roupa%2Bde%2Bver%C3%A3o%2BGatito.png"
    10 = "https://github.com/rosanamassahud/cursoPowerBI_base1/blob/main/img_petshop/Coleira%2BGatito.png"
    11 = "https://github.com/rosanamassahud/cursoPowerBI_base1/blob/main/img_petshop/Guia%2BGatito.png"
    12 = "https://github.com/rosanamassahud/cursoPowerBI_base1/blob/main/img_petshop/Ra%C3%A7%C3%A3o%2Bfilhote%2BCachorrito.png"
    13 = "https://github.com/rosanamassahud/cursoPowerBI_base1/blob/main/img_petshop/Ra%C3%A7%C3%A3o%2Badulto%2BCachorrito.png"
    14 = "https://github.com/rosanamassahud/cursoPowerBI_base1/blob/main/img_petshop/Ra%C3%A7%C3%A3o%2Bsenior%2BCachorrito.png"
    15 = "https://github.com/rosanamassahud/cursoPowerBI_base1/blob/main/img_petshop/Petisco%2BCachorrito.png"
    16 = "https://github.com/rosanamassahud/cursoPowerBI_base1/blob/main/img_petshop/Bola%2BCachorrito.png"
    17 = "https://github.com/rosanamassahud/cursoPowerBI_base1/blob/main/img_petshop/Osso%2BCachorrito.png"
    18 = "https://github.com/rosanamassahud/cursoPowerBI_base1/blob/main/img_petshop/Roupa%2Bde%2Binverno%2BCachorrito.png"
    19 = "https://github.com/rosanamassahud/cursoPowerBI_base1/blob/main/img_petshop/Roupa%2Bde%2Bver%C3%A3o%2BCachorrito.png"
    20 = "https://github.com/rosanamassahud/cursoPowerBI_base1/blob/main/img_petshop/Coleira%2BCachorrito.png"
    21 = "https://github.com/rosanamassahud/cursoPowerBI_base1/blob/main/img_petshop/Guia%2BCachorrito.png"
}

# Row 3 keeps the "Hiperlink" visual style (underlined, themed blue text)
# but - unlike every other row - does NOT get an actual clickable
# hyperlink relationship, matching the source workbook exactly.
$noLinkRows = @(3)

foreach ($row in 2..21) {
    $cell = $ws.Cells.Item($row, 6)
    $url = $urls[$row]
    $cell.Value = $url
    if ($noLinkRows -contains $row) {
        $cell.Style = "Hyperlink"
    } else {
        $ws.Hyperlinks.Add($cell, $url) | Out-Null
    }
}

# Restore the selection left behind on F5 by the author.
$ws.Range("F5").Select() | Out-Null
